$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- Row 1: turn the old "duplicate data" header row into real column headers ---
# (matching the layout used by the 土地/建物/股票 sheets: name, capacity, owner,
# register_date, register_reason, acquire_value, property_category, category,
# date, legislator_name, legislator_id, source_file, index)
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"

# New header cells H1:N1 - first copy the style from an existing header cell (B1),
# then fill in the text.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("H1:N1").PasteSpecial(-4122) | Out-Null

$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Row 2 (BluebirdSylphyGllTS(客車)) ---
$ws.Range("B2").Value = "BluebirdSylphyGllTS(客車)"
$ws.Range("E2").Value = "98年10月12日"

$ws.Range("B2").Copy() | Out-Null
$ws.Range("H2:N2").PasteSpecial(-4122) | Out-Null

$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"
# force text so "2013-12-24" isn't reinterpreted as a date serial number
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2013-12-24"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("J2").PasteSpecial(-4122) | Out-Null
$ws.Range("K2").Value = "林國正"
$ws.Range("L2").Value = 1742
$ws.Range("M2").Value = "tmp399c1"
$ws.Range("N2").Value = 29

# --- Row 3 (納智捷L91MLD) ---
$ws.Range("B3").Value = "納智捷L91MLD"
$ws.Range("D3").Value = "林國正"
$ws.Range("E3").Value = "102年01月15曰"
$ws.Range("F3").Value = "買賣."

$ws.Range("B3").Copy() | Out-Null
$ws.Range("H3:N3").PasteSpecial(-4122) | Out-Null

$ws.Range("H3").Value = "land"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").NumberFormat = "@"
$ws.Range("J3").Value = "2013-12-24"
$ws.Range("B3").Copy() | Out-Null
$ws.Range("J3").PasteSpecial(-4122) | Out-Null
$ws.Range("K3").Value = "林國正"
$ws.Range("L3").Value = 1742
$ws.Range("M3").Value = "tmp399c1"
$ws.Range("N3").Value = 30
